# single_function_8params.xlsx - update eval data of LLM
#
# The "detail" column (D2:D51) on Sheet1 all share one de-duplicated shared
# string that records token/function-call accounting for the eval run.
# The "completionTokens" field was dropped from that JSON blob, so every
# cell that pointed at the old shared string needs to be rewritten with the
# trimmed JSON (rewriting the whole range in one shot keeps Excel's shared
# string de-duplication working instead of forking a second <si> entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDetail = '{"fcCount":1,"fcInfo":{"apiair-conditionerupdate_POST":8}}'
$ws.Range("D2:D51").Value = $newDetail

# Reflect the editor's new viewport/selection on Sheet1: the cursor moved
# from D2 to a block selection starting at E2.
[void]$ws.Range("E2:BB86").Select()
